$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(308051846, Eyal  Sofer: -1,2)"
$ws.Range("B1").Value = "(312049950, Molham  Peretz: -7,-9)"
$ws.Range("C1").Value = "(308073899, Anan  Kirshenbaum: -8,0)"
$ws.Range("D1").Value = "(318869187, Soaad  Leibovich: 3,-3)"
$ws.Range("E1").Value = "(205898513, Asaf  Braymok: 2,3)"
$ws.Range("F1").Value = "(318428158, Tal  Asulin: 6,9)"
$ws.Range("G1").Value = "(316028364, Sami  Castro: 5,4)"

$ws.Range("A3").Value = "cost: 338.6127685688407"
$ws.Range("A4").Value = "time: 62.72255371376816"
